$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 37; existing rows 37..139 shift down to 38..140.
$ws.Rows.Item(37).Insert()

# Populate the newly inserted row 37 with the new record's data.
$fecha = Get-Date -Year 2022 -Month 4 -Day 21 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(37, 1).Value = 8
$ws.Cells.Item(37, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(37, 3).Value = "Coquimbo"
$ws.Cells.Item(37, 4).Value = $fecha.Date
$ws.Cells.Item(37, 5).Value = 4
$ws.Cells.Item(37, 6).Value = 100112044
$ws.Cells.Item(37, 7).Value = "Perejil"
$ws.Cells.Item(37, 8).Value = "Sin especificar"
$ws.Cells.Item(37, 9).Value = "Primera"
$ws.Cells.Item(37, 10).Value = 2400
$ws.Cells.Item(37, 11).Value = 2000
$ws.Cells.Item(37, 12).Value = 2500
$ws.Cells.Item(37, 13).Value = 2250
$ws.Cells.Item(37, 14).Value = "`$/atado 1 a 1,5 kilos"
$ws.Cells.Item(37, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(37, 16).Value = 1500
$ws.Cells.Item(37, 17).Value = 1.5
$ws.Cells.Item(37, 18).Value = "Hortaliza"
